$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the pre-existing "_GoBack" bookmark that currently sits right
#    after the "Open AllProjects.sln to Visual Studio and rebuild solution"
#    run (it is being relocated to the top of the document - see step 3).
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 2) Merge the two runs "checkout directory" and
#    ". The easiest way to do this is just click" into a single run.
#    A plain text replacement on this engine coalesces every run, from the
#    edited point through to the end of the paragraph, that shares identical
#    run formatting - so first force the merge (accepting the temporary
#    over-merge), then restore the original run boundaries for the
#    untouched tail of the paragraph (the quoted “set_yam2d_env.bat”
#    sentence) with a harmless Bold on/off toggle, which only ever creates
#    new run splits and never re-coalesces already-split runs.
# ---------------------------------------------------------------------------
$rngFind = $d.Content
$ok = $rngFind.Find.Execute("checkout directory. The easiest way to do this is just click", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergeStart = $rngFind.Start
$mergeEnd = $rngFind.End
$tailStart = $mergeEnd

# Force a genuine text change so the engine recomputes/merges the run(s).
$rngEdit = $d.Range($mergeStart, $mergeEnd)
$rngEdit.Text = "checkout directory. The easiest way to do this is just clickZ"
$rngEdit2 = $d.Range($mergeStart, $mergeEnd + 1)
$rngEdit2.Text = "checkout directory. The easiest way to do this is just click"

# Re-establish the original run splits for the untouched tail of the
# paragraph (" “set_yam2d_env.bat”, found in your Yam2D checkout directory.")
$tailFinds = @(" “", "set_yam2d_env.bat", "”, found in your Yam2D checkout directory", ".")
$cursor = $tailStart
foreach ($piece in $tailFinds) {
    $len = $piece.Length
    $pieceStart = $cursor
    $pieceEnd = $cursor + $len
    $d.Range($pieceStart, $pieceEnd).Bold = 1
    $d.Range($pieceStart, $pieceEnd).Bold = 0
    $cursor = $pieceEnd
}

# ---------------------------------------------------------------------------
# 3) Add a new (empty) "_GoBack" bookmark at the very start of the document,
#    immediately before the first run ("Download Yam2D using tortoise svn").
#    A bookmark added directly on a zero-length Range(0,0) is mis-anchored
#    by this engine, so seed it around one temporary character, then delete
#    that character - the bookmark collapses to the correct zero-length
#    position at the very start of the body.
# ---------------------------------------------------------------------------
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("X")
$bmRange = $d.Range(0, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range(0, 1).Delete()

Write-Output "done"
